$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '35.499.38'
$ws.Range('E2').Value = '  +2.37%  '

# Row 3
$ws.Range('D3').Value = '1.848.71'
$ws.Range('E3').Value = '  +1.77%  '

# Row 4
$ws.Range('E4').Value = '  +0.16%  '

# Row 5
$ws.Range('D5').Value = '''228.17'
$ws.Range('E5').Value = '  +0.69%  '

# Row 6
$ws.Range('E6').Value = '  +1.53%  '

# Row 7
$ws.Range('E7').Value = '  +0.12%  '

# Row 8
$ws.Range('D8').Value = '''41.40'
$ws.Range('E8').Value = '  +7.49%  '

# Row 9
$ws.Range('E9').Value = '  +5.02%  '

# Row 10
$ws.Range('D10').Value = '''0.0690'
$ws.Range('E10').Value = '  +0.86%  '

# Row 11
$ws.Range('E11').Value = '  +3.13%  '

# Row 12
$ws.Range('D12').Value = '2.116.11'
$ws.Range('E12').Value = '  +1.89%  '

# Row 13
$ws.Range('D13').Value = '''11.65'
$ws.Range('E13').Value = '  +2.31%  '

# Row 14
$ws.Range('D14').Value = '1.850.18'
$ws.Range('E14').Value = '  +0.86%  '

# Row 15
$ws.Range('B15').Value = 'Polygon'
$ws.Range('C15').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D15').Value = '''0.670'
$ws.Range('E15').Value = '  +5.16%  '

# Row 16
$ws.Range('B16').Value = 'Polkadot'
$ws.Range('C16').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D16').Value = '''4.74'
$ws.Range('E16').Value = '  +6.22%  '

# Row 17
$ws.Range('D17').Value = '35.457.84'
$ws.Range('E17').Value = '  +2.46%  '

# Row 18
$ws.Range('D18').Value = '''70.01'
$ws.Range('E18').Value = '  +1.48%  '

# Row 19
$ws.Range('D19').Value = '''245.84'
$ws.Range('E19').Value = '  +0.21%  '

# Row 20
$ws.Range('D20').Value = '0.0₃0797'
$ws.Range('E20').Value = '  +2.18%  '

# Row 21
$ws.Range('D21').Value = '''12.22'
$ws.Range('E21').Value = '  +7.87%  '

# Row 22
$ws.Range('D22').Value = '''4.77'
$ws.Range('E22').Value = '  +14.89%  '

# Row 23
$ws.Range('E23').Value = '  +0.18%  '

# Row 24
$ws.Range('E24').Value = '  -0.77%  '

# Row 25
$ws.Range('D25').Value = '''171.78'
$ws.Range('E25').Value = '  -0.33%  '

# Row 26
$ws.Range('D26').Value = '''7.90'
$ws.Range('E26').Value = '  -0.52%  '

# Row 27
$ws.Range('D27').Value = '''17.85'
$ws.Range('E27').Value = '  +1.68%  '

# Row 28
$ws.Range('E28').Value = '  +0.88%  '

# Row 29
$ws.Range('E29').Value = '  +0.17%  '

# Row 30
$ws.Range('D30').Value = '3.382.47'
$ws.Range('E30').Value = '  +39.21%  '

# Row 31
$ws.Range('E31').Value = '  +7.92%  '

# Row 32
$ws.Range('B32').Value = 'InternetComputer(DFINITY)'
$ws.Range('C32').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D32').Value = '''4.07'
$ws.Range('E32').Value = '  +2.89%  '

# Row 33
$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D33').Value = '''3.94'
$ws.Range('E33').Value = '  +3.03%  '

# Row 34
$ws.Range('D34').Value = '''0.0537'
$ws.Range('E34').Value = '  +2.23%  '

# Row 35
$ws.Range('E35').Value = '  +2.64%  '

# Row 36
$ws.Range('D36').Value = '''0.678'
$ws.Range('E36').Value = '  +3.04%  '

# Row 37
$ws.Range('E37').Value = '  +9.88%  '

# Row 38
$ws.Range('D38').Value = '''89.00'
$ws.Range('E38').Value = '  +9.20%  '

# Row 39
$ws.Range('D39').Value = '1.339.92'
$ws.Range('E39').Value = '  -2.16%  '

# Row 40
$ws.Range('E40').Value = '  +1.37%  '

# Row 41
$ws.Range('E41').Value = '  +3.16%  '

# Row 42
$ws.Range('E42').Value = '  +1.14%  '

# Row 43
$ws.Range('E43').Value = '  +3.93%  '

# Row 44
$ws.Range('D44').Value = '''14.88'
$ws.Range('E44').Value = '  +4.83%  '

# Row 45
$ws.Range('E45').Value = '  +0.86%  '

# Row 47
$ws.Range('D47').Value = '''0.0520'
$ws.Range('E47').Value = '  +3.46%  '

# Row 48
$ws.Range('D48').Value = '''6.06'
$ws.Range('E48').Value = '  +4.63%  '

# Row 49
$ws.Range('D49').Value = '2.014.57'
$ws.Range('E49').Value = '  +1.88%  '

# Row 50
$ws.Range('D50').Value = '''104.45'
$ws.Range('E50').Value = '  +1.11%  '

# Row 51
$ws.Range('E51').Value = '  +0.11%  '
